$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# 1) Swap the full data (columns B:AC) of row 2 and row 3.
#    (column A, the running "id", stays 0/1 in order)
# -----------------------------------------------------------------------
$row2 = $ws.Range("B2:AC2").Value2
$row3 = $ws.Range("B3:AC3").Value2

$ws.Range("B2:AC2").Value2 = $row3
$ws.Range("B3:AC3").Value2 = $row2

# -----------------------------------------------------------------------
# 2) Insert a brand new match as row 153 (pushing the former rows
#    153-155 down to 154-156).
# -----------------------------------------------------------------------
# xlShiftDown = -4121
$ws.Rows(153).Insert(-4121) | Out-Null

# Re-apply the same formatting used by every other row: bold / centered /
# bordered "id" cell in column A, and the custom date format in column E.
$ws.Range("A153").Font.Bold = $true
$ws.Range("A153").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A153").VerticalAlignment = -4160    # xlTop
$ws.Range("A153").Borders.LineStyle = 1
$ws.Range("E153").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A153").Value2 = 151
$ws.Range("B153").Value2 = 8100619
$ws.Range("C153").Value2 = "Slovenia Prva Liga"
$ws.Range("D153").Value2 = "Slovenia Prva Liga"
$ws.Range("E153").Value2 = 45400.45833333334
$ws.Range("F153").Value2 = "NK Domzale"
$ws.Range("G153").Value2 = "FC Koper"
$ws.Range("H153").Value2 = 1
$ws.Range("I153").Value2 = 0
$ws.Range("J153").Value2 = "H"
$ws.Range("K153").Value2 = 3.1
$ws.Range("L153").Value2 = 3.5
$ws.Range("M153").Value2 = 2.05
$ws.Range("N153").Value2 = 3
$ws.Range("O153").Value2 = 3.5
$ws.Range("P153").Value2 = 2.15
$ws.Range("Q153").Value2 = 0.25
$ws.Range("R153").Value2 = 1.9
$ws.Range("S153").Value2 = 1.9
$ws.Range("T153").Value2 = 2.75
$ws.Range("U153").Value2 = 1.925
$ws.Range("V153").Value2 = 1.875
$ws.Range("W153").Value2 = 2
$ws.Range("X153").Value2 = -1
$ws.Range("Y153").Value2 = -1
$ws.Range("Z153").Value2 = 0.8999999999999999
$ws.Range("AA153").Value2 = -1
$ws.Range("AB153").Value2 = -1
$ws.Range("AC153").Value2 = 0.875

# The running "id" in column A (= row number - 2) must be bumped by one
# for every row that got pushed down by the insertion above.
$ws.Range("A154").Value2 = 152
$ws.Range("A155").Value2 = 153
$ws.Range("A156").Value2 = 154
